# Planificacion.xlsx - "Comenzada idea de spawning en clusters, falta codigo"
# Insert a new column (H) on the "Enemigos" sheet holding the average damage
# dealt by the player to the enemy (2 for every enemy for now), shifting the
# old "VARIABLES JUGADOR" block one column to the right (I:J -> J:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enemigos")

# Shift the "VARIABLES JUGADOR" block (old col I/J) one column to the right
# so there is room for the new "Daño Medio (de jugador)" column at H.
$ws.Columns("I:I").Insert()

# New column header + values.
$ws.Range("H2").Value = "Daño Medio (de jugador)"
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 2
$ws.Range("H5").Value = 2

# Match the centered style used by the other header cells / the merged title.
$ws.Range("H1").HorizontalAlignment = -4108

# Give the new column a sensible best-fit width like its neighbours
# (closest the engine's pixel grid can reach to the authored 23.28515625).
$ws.Columns("H:H").ColumnWidth = 22.5

# Move the active selection like in the authored workbook.
$ws.Range("H7").Select()

$wb.Save()
